# Scheduled market-price refresh: updates currentAveragePrice(NQ/HQ) and
# LevePrice/LeveProfit columns (H:N) across all job sheets with freshly
# scraped values. Purely literal numeric overwrites, no formulas involved.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 316.26315
$ws.Range("I2").Value = 197.64285
$ws.Range("J2").Value = 648.4
$ws.Range("K2").Value = 197.64285
$ws.Range("L2").Value = 648.4
$ws.Range("M2").Value = -84.64285000000001
$ws.Range("N2").Value = -874.4
$ws.Range("H31").Value = 1980.8182
$ws.Range("I31").Value = 214.83333
$ws.Range("J31").Value = 4100
$ws.Range("K31").Value = 644.49999
$ws.Range("L31").Value = 12300
$ws.Range("M31").Value = -414.49999
$ws.Range("N31").Value = -12760
$ws.Range("H40").Value = 4540.9565
$ws.Range("I40").Value = 3662.9285
$ws.Range("J40").Value = 5906.778
$ws.Range("K40").Value = 3662.9285
$ws.Range("L40").Value = 5906.778
$ws.Range("M40").Value = -3487.9285
$ws.Range("N40").Value = -6256.778
$ws.Range("H51").Value = 3435
$ws.Range("I51").Value = 2850
$ws.Range("J51").Value = 3581.25
$ws.Range("K51").Value = 2850
$ws.Range("L51").Value = 3581.25
$ws.Range("M51").Value = -2366
$ws.Range("N51").Value = -4549.25
$ws.Range("H80").Value = 1042.2354
$ws.Range("I80").Value = 696.125
$ws.Range("K80").Value = 2088.375
$ws.Range("M80").Value = -1090.375
$ws.Range("H83").Value = 1042.2354
$ws.Range("I83").Value = 696.125
$ws.Range("K83").Value = 6265.125
$ws.Range("M83").Value = -1273.125
$ws.Range("H86").Value = 4285.2856
$ws.Range("I86").Value = 3833.3333
$ws.Range("J86").Value = 4624.25
$ws.Range("K86").Value = 3833.3333
$ws.Range("L86").Value = 4624.25
$ws.Range("M86").Value = -2710.3333
$ws.Range("N86").Value = -6870.25
$ws.Range("H88").Value = 1009
$ws.Range("I88").Value = 1950
$ws.Range("J88").Value = 632.6
$ws.Range("K88").Value = 1950
$ws.Range("L88").Value = 632.6
$ws.Range("M88").Value = -1544
$ws.Range("N88").Value = -1444.6
$ws.Range("H89").Value = 4285.2856
$ws.Range("I89").Value = 3833.3333
$ws.Range("J89").Value = 4624.25
$ws.Range("K89").Value = 19166.6665
$ws.Range("L89").Value = 23121.25
$ws.Range("M89").Value = -13550.6665
$ws.Range("N89").Value = -34353.25
$ws.Range("H91").Value = 1009
$ws.Range("I91").Value = 1950
$ws.Range("J91").Value = 632.6
$ws.Range("K91").Value = 1950
$ws.Range("L91").Value = 632.6
$ws.Range("M91").Value = -546
$ws.Range("N91").Value = -3440.6
$ws.Range("H94").Value = 6669.0713
$ws.Range("I94").Value = 6669.0713
$ws.Range("K94").Value = 6669.0713
$ws.Range("M94").Value = -6218.0713
$ws.Range("H117").Value = 125000
$ws.Range("J117").Value = 125000
$ws.Range("L117").Value = 125000
$ws.Range("N117").Value = -134178
$ws.Range("H137").Value = 3412.3809
$ws.Range("I137").Value = 2127.4443
$ws.Range("J137").Value = 4376.0835
$ws.Range("K137").Value = 6382.3329
$ws.Range("L137").Value = 13128.2505
$ws.Range("M137").Value = -3832.3329
$ws.Range("N137").Value = -18228.2505

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4018.439
$ws.Range("I32").Value = 4018.439
$ws.Range("K32").Value = 4018.439
$ws.Range("M32").Value = -3731.439
$ws.Range("H88").Value = 3400
$ws.Range("I88").Value = 4750
$ws.Range("J88").Value = 2500
$ws.Range("K88").Value = 4750
$ws.Range("L88").Value = 2500
$ws.Range("M88").Value = -4344
$ws.Range("N88").Value = -3312
$ws.Range("H91").Value = 3400
$ws.Range("I91").Value = 4750
$ws.Range("J91").Value = 2500
$ws.Range("K91").Value = 4750
$ws.Range("L91").Value = 2500
$ws.Range("M91").Value = -3346
$ws.Range("N91").Value = -5308
$ws.Range("H97").Value = 576.7
$ws.Range("I97").Value = 554.4211
$ws.Range("K97").Value = 554.4211
$ws.Range("M97").Value = -58.42110000000002
$ws.Range("H110").Value = 2996.6155
$ws.Range("I110").Value = 1369.875
$ws.Range("K110").Value = 1369.875
$ws.Range("M110").Value = 675.125
$ws.Range("H132").Value = 3762.8462
$ws.Range("I132").Value = 1799.6
$ws.Range("K132").Value = 5398.799999999999
$ws.Range("M132").Value = -2868.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H100").Value = 13850
$ws.Range("J100").Value = 13850
$ws.Range("L100").Value = 13850
$ws.Range("N100").Value = -16014

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.652172
$ws.Range("J2").Value = 35.77778
$ws.Range("L2").Value = 214.66668
$ws.Range("N2").Value = -440.66668
$ws.Range("H109").Value = 127080.125
$ws.Range("J109").Value = 2196.6667
$ws.Range("L109").Value = 6590.000100000001
$ws.Range("N109").Value = -8670.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2686.75
$ws.Range("I80").Value = 2499.4285
$ws.Range("K80").Value = 2499.4285
$ws.Range("M80").Value = -1501.4285
$ws.Range("H83").Value = 2686.75
$ws.Range("I83").Value = 2499.4285
$ws.Range("K83").Value = 12497.1425
$ws.Range("M83").Value = -7505.1425
$ws.Range("H92").Value = 7050
$ws.Range("I92").Value = 1000
$ws.Range("J92").Value = 8562.5
$ws.Range("K92").Value = 1000
$ws.Range("L92").Value = 8562.5
$ws.Range("M92").Value = 872
$ws.Range("N92").Value = -12306.5
$ws.Range("H132").Value = 110713.9
$ws.Range("I132").Value = 135392.5
$ws.Range("K132").Value = 406177.5
$ws.Range("M132").Value = -403647.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 9990
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H33").Value = 9646.75
$ws.Range("J33").Value = 8617
$ws.Range("L33").Value = 8617
$ws.Range("N33").Value = -9197
$ws.Range("H43").Value = 5975
$ws.Range("J43").Value = 5975
$ws.Range("L43").Value = 5975
$ws.Range("N43").Value = -6361
$ws.Range("H82").Value = 4382.8237
$ws.Range("I82").Value = 3804.5715
$ws.Range("K82").Value = 3804.5715
$ws.Range("M82").Value = -3443.5715
$ws.Range("H85").Value = 4382.8237
$ws.Range("I85").Value = 3804.5715
$ws.Range("K85").Value = 3804.5715
$ws.Range("M85").Value = -2556.5715
$ws.Range("H93").Value = 2259.5454
$ws.Range("I93").Value = 2259.5454
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2259.5454
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1011.5454
$ws.Range("N93").ClearContents()
$ws.Range("H136").Value = 4950.9165
$ws.Range("I136").Value = 4741.2
$ws.Range("K136").Value = 14223.6
$ws.Range("M136").Value = -11673.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 549.2
$ws.Range("I107").Value = 515.75
$ws.Range("K107").Value = 1547.25
$ws.Range("M107").Value = 372.75
$ws.Range("H136").Value = 3833.6667
$ws.Range("I136").Value = 2902.5
$ws.Range("J136").Value = 5696
$ws.Range("K136").Value = 8707.5
$ws.Range("L136").Value = 17088
$ws.Range("M136").Value = -6157.5
$ws.Range("N136").Value = -22188
